$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A7").Value = "oo0871"
$ws.Range("C7").Value = "매니저"
